$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 317.2
$ws.Range("I4").Value = 285
$ws.Range("K4").Value = 285
$ws.Range("M4").Value = -171
$ws.Range("H17").Value = 2225.3076
$ws.Range("J17").Value = 2225.3076
$ws.Range("L17").Value = 6675.9228
$ws.Range("N17").Value = -7011.9228
$ws.Range("H33").Value = 222.07408
$ws.Range("J33").Value = 1399.5
$ws.Range("L33").Value = 1399.5
$ws.Range("N33").Value = -1857.5
$ws.Range("H43").Value = 1511.3684
$ws.Range("J43").Value = 1672.7693
$ws.Range("L43").Value = 1672.7693
$ws.Range("N43").Value = -1810.7693
$ws.Range("H132").Value = 372114.22
$ws.Range("I132").Value = 1641.591
$ws.Range("K132").Value = 4924.772999999999
$ws.Range("M132").Value = -2394.772999999999
$ws.Range("H138").Value = 4414.049
$ws.Range("I138").Value = 3127.8
$ws.Range("K138").Value = 9383.400000000001
$ws.Range("M138").Value = -4243.400000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15754.678
$ws.Range("I32").Value = 5504.913
$ws.Range("K32").Value = 5504.913
$ws.Range("M32").Value = -5217.913
$ws.Range("H42").Value = 25049.182
$ws.Range("I42").Value = 16319.2
$ws.Range("J42").Value = 32324.166
$ws.Range("K42").Value = 16319.2
$ws.Range("L42").Value = 32324.166
$ws.Range("M42").Value = -15833.2
$ws.Range("N42").Value = -33296.166
$ws.Range("H45").Value = 1715.4286
$ws.Range("I45").Value = 1768
$ws.Range("K45").Value = 1768
$ws.Range("M45").Value = -1391
$ws.Range("H61").Value = 2793.7585
$ws.Range("I61").Value = 2206.5715
$ws.Range("J61").Value = 4335.125
$ws.Range("K61").Value = 2206.5715
$ws.Range("L61").Value = 4335.125
$ws.Range("M61").Value = -1994.5715
$ws.Range("N61").Value = -4759.125
$ws.Range("H74").Value = 831.381
$ws.Range("I74").Value = 813.6316
$ws.Range("K74").Value = 813.6316
$ws.Range("M74").Value = 60.36839999999995
$ws.Range("H77").Value = 831.381
$ws.Range("I77").Value = 813.6316
$ws.Range("K77").Value = 4068.158
$ws.Range("M77").Value = 299.8419999999996
$ws.Range("H136").Value = 2793.7585
$ws.Range("I136").Value = 2206.5715
$ws.Range("J136").Value = 4335.125
$ws.Range("K136").Value = 6619.7145
$ws.Range("L136").Value = 13005.375
$ws.Range("M136").Value = -4069.7145
$ws.Range("N136").Value = -18105.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1476
$ws.Range("I64").Value = 1476
$ws.Range("K64").Value = 1476
$ws.Range("M64").Value = -1251
$ws.Range("H67").Value = 1476
$ws.Range("I67").Value = 1476
$ws.Range("K67").Value = 1476
$ws.Range("M67").Value = -696
$ws.Range("H86").Value = 2035.6666
$ws.Range("I86").Value = 1600
$ws.Range("J86").Value = 2471.3333
$ws.Range("K86").Value = 1600
$ws.Range("L86").Value = 2471.3333
$ws.Range("M86").Value = -477
$ws.Range("N86").Value = -4717.3333
$ws.Range("H89").Value = 2035.6666
$ws.Range("I89").Value = 1600
$ws.Range("J89").Value = 2471.3333
$ws.Range("K89").Value = 8000
$ws.Range("L89").Value = 12356.6665
$ws.Range("M89").Value = -2384
$ws.Range("N89").Value = -23588.6665
$ws.Range("H105").Value = 3591.7407
$ws.Range("I105").Value = 2972
$ws.Range("J105").Value = 4087.5334
$ws.Range("K105").Value = 2972
$ws.Range("L105").Value = 4087.5334
$ws.Range("M105").Value = -1225
$ws.Range("N105").Value = -7581.5334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 20007
$ws.Range("I17").Value = 20007
$ws.Range("K17").Value = 20007
$ws.Range("M17").Value = -19833
$ws.Range("H22").Value = 717.3333
$ws.Range("I22").Value = 581.6667
$ws.Range("K22").Value = 581.6667
$ws.Range("M22").Value = -231.6667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2293021
$ws.Range("J4").Value = 499.33334
$ws.Range("L4").Value = 1498.00002
$ws.Range("N4").Value = -1722.00002
$ws.Range("H11").Value = 1263
$ws.Range("I11").Value = 968.8182
$ws.Range("K11").Value = 2906.4546
$ws.Range("M11").Value = -2766.4546
$ws.Range("H129").Value = 7263.364
$ws.Range("I129").Value = 2166.3333
$ws.Range("J129").Value = 9174.75
$ws.Range("K129").Value = 6498.999899999999
$ws.Range("L129").Value = 27524.25
$ws.Range("M129").Value = -1498.999899999999
$ws.Range("N129").Value = -37524.25
$ws.Range("H137").Value = 1999.5
$ws.Range("J137").Value = 1999.5
$ws.Range("L137").Value = 5998.5
$ws.Range("N137").Value = -16198.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4397.4707
$ws.Range("I132").Value = 4072.6667
$ws.Range("J132").Value = 5177
$ws.Range("K132").Value = 12218.0001
$ws.Range("L132").Value = 15531
$ws.Range("M132").Value = -9688.000100000001
$ws.Range("N132").Value = -20591

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5825.6313
$ws.Range("I7").Value = 6054.091
$ws.Range("K7").Value = 6054.091
$ws.Range("M7").Value = -5942.091
$ws.Range("H40").Value = 7129.4707
$ws.Range("I40").Value = 3689.111
$ws.Range("J40").Value = 10999.875
$ws.Range("K40").Value = 3689.111
$ws.Range("L40").Value = 10999.875
$ws.Range("M40").Value = -3553.111
$ws.Range("N40").Value = -11271.875
$ws.Range("H96").Value = 70197
$ws.Range("J96").Value = 70197
$ws.Range("L96").Value = 70197
$ws.Range("N96").Value = -75689
$ws.Range("H126").Value = 5825.6313
$ws.Range("I126").Value = 6054.091
$ws.Range("K126").Value = 18162.273
$ws.Range("M126").Value = -15692.273
$ws.Range("H136").Value = 3614.2307
$ws.Range("I136").Value = 3212.8572
$ws.Range("K136").Value = 9638.571599999999
$ws.Range("M136").Value = -7088.571599999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 13394
$ws.Range("I62").Value = 7591.5
$ws.Range("K62").Value = 7591.5
$ws.Range("M62").Value = -6967.5
$ws.Range("H65").Value = 13394
$ws.Range("I65").Value = 7591.5
$ws.Range("K65").Value = 37957.5
$ws.Range("M65").Value = -34837.5
$ws.Range("H96").Value = 3838.3
$ws.Range("I96").Value = 3838.3
$ws.Range("K96").Value = 3838.3
$ws.Range("M96").Value = -2465.3
$ws.Range("H132").Value = 1291.6571
$ws.Range("I132").Value = 1291.6571
$ws.Range("K132").Value = 3874.9713
$ws.Range("M132").Value = -1344.9713
